$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells in row 8 for the columns that survive the column delete below,
# BEFORE deleting, since the delete will shift column references.
$ws.Range("D8").Value = "Product S/N"
$ws.Range("E8").Value = "PCBA S/N"
$ws.Range("F8").Value = "Device ID"
$ws.Range("G8").Value = "WO No."

# Widen column E (now "PCBA S/N") and column G (now "WO No.") - they no longer
# auto-size ("best fit") to their header text, they get an explicit width instead.
$ws.Columns.Item(5).ColumnWidth = 23.6666666666667
$ws.Columns.Item(7).ColumnWidth = 13.8333333333333

# Delete column H ("Defect Desc.") entirely - shifts old I (Created By) -> H,
# old J (Date Created) -> I, matching the template's new 8-column layout.
$ws.Columns.Item(8).Delete()

# Re-point the selection like the edited workbook.
$ws.Range("H18").Select() | Out-Null
